$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 110; existing rows 110-152 shift down to 111-153,
# carrying their formatting (including the date style on column D) with them.
$ws.Rows("110:110").Insert()

# Populate the newly inserted row 110 with the new weekly record.
$ws.Cells.Item(110, 1).Value = 5
$ws.Cells.Item(110, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(110, 3).Value = "Maule"
$ws.Cells.Item(110, 4).Value = 44917
$ws.Cells.Item(110, 5).Value = 7
$ws.Cells.Item(110, 6).Value = 100112030
$ws.Cells.Item(110, 7).Value = "Poroto granado"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 150
$ws.Cells.Item(110, 11).Value = 32000
$ws.Cells.Item(110, 12).Value = 32000
$ws.Cells.Item(110, 13).Value = 32000
$ws.Cells.Item(110, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(110, 15).Value = "Región del Maule"
$ws.Cells.Item(110, 16).Value = 1280
$ws.Cells.Item(110, 17).Value = 25
$ws.Cells.Item(110, 18).Value = "Hortaliza"
